$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers for the data columns that carry this week-to-week data.
# D=4 (Fecha), H=8 (Variedad), I=9 (Calidad), J=10 (Volumen), K=11 (Precio minimo),
# L=12 (Precio maximo), M=13 (Precio promedio ponderado), N=14 (Unidad de comercializacion),
# O=15 (Origen), P=16 (Precio $/Kg)
$cols = @(4,8,9,10,11,12,13,14,15,16)

# Snapshot every existing data row (41..175) BEFORE any writes, since the
# transform shifts each row's data down into the next row.
$oldData = @{}
for ($r = 41; $r -le 175; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $oldData[$r] = $rowVals
}

# Row 41 keeps its own Volumen/Precio/Calidad/Variedad/Origen values, but gets
# a brand-new date and its "Unidad de comercializacion" flips to "$/caja 10 kilos".
$ws.Cells.Item(41, 4).Value = 44487
$ws.Cells.Item(41, 14).Value = "$/caja 10 kilos"

# Rows 42..175 each take on what used to be the row above them (41..174).
for ($r = 42; $r -le 175; $r++) {
    $src = $oldData[$r - 1]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $src[$c]
    }
}

# New row 176 takes the constant columns from row 175 (they're the same across
# the whole sheet) plus the data that used to live in row 175.
$ws.Cells.Item(176, 1).Value = $ws.Cells.Item(175, 1).Value2
$ws.Cells.Item(176, 2).Value = $ws.Cells.Item(175, 2).Value2
$ws.Cells.Item(176, 3).Value = $ws.Cells.Item(175, 3).Value2
$ws.Cells.Item(176, 5).Value = $ws.Cells.Item(175, 5).Value2
$ws.Cells.Item(176, 6).Value = $ws.Cells.Item(175, 6).Value2
$ws.Cells.Item(176, 7).Value = $ws.Cells.Item(175, 7).Value2
$ws.Cells.Item(176, 17).Value = $ws.Cells.Item(175, 17).Value2
$ws.Cells.Item(176, 18).Value = $ws.Cells.Item(175, 18).Value2

$srcLast = $oldData[175]
foreach ($c in $cols) {
    $ws.Cells.Item(176, $c).Value = $srcLast[$c]
}

# The "Fecha" column carries a date number format on every data row; copy it
# onto the freshly created row 176 cell so it matches its neighbours.
$ws.Cells.Item(176, 4).NumberFormat = $ws.Cells.Item(175, 4).NumberFormat
